$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Rows 3-8: content unchanged; only the new "Date Completed" value
# (col I) is populated, reusing the date style already used by I15.
# ---------------------------------------------------------------
$ws.Range("I15").Copy()
$ws.Range("I3:I8").PasteSpecial(-4122)
$ws.Cells.Item(3,9).Value = 43860
$ws.Cells.Item(4,9).Value = 43860
$ws.Cells.Item(5,9).Value = 43860
$ws.Cells.Item(6,9).Value = 43861
$ws.Cells.Item(7,9).Value = 43861
$ws.Cells.Item(8,9).Value = 43861

# ---------------------------------------------------------------
# New rows 23-24: clone formatting from row 17 (B:I) first so the
# freshly-created cells carry the same styles (s="1" / s="2") as
# the rest of the table, then fill in their values below.
# ---------------------------------------------------------------
$ws.Range("B17:I17").Copy()
$ws.Range("B23:I24").PasteSpecial(-4122)

# ---------------------------------------------------------------
# Rows 9-24: scenario table re-sorted into Watershed/Species order;
# rewrite B:H for each row from the new layout.
# ---------------------------------------------------------------
# Row 9: scenario 7
$ws.Cells.Item(9,2).Value = 7
$ws.Cells.Item(9,3).Value = "Pahsimeroi"
$ws.Cells.Item(9,4).Value = "Chinook"
$ws.Cells.Item(9,5).Value = "Juvenile"
$ws.Cells.Item(9,6).Value = "Summer"
$ws.Cells.Item(9,7).Clear()
$ws.Cells.Item(9,8).Clear()

# Row 10: scenario 8
$ws.Cells.Item(10,2).Value = 8
$ws.Cells.Item(10,3).Value = "Pahsimeroi"
$ws.Cells.Item(10,4).Value = "Chinook"
$ws.Cells.Item(10,5).Value = "Juvenile"
$ws.Cells.Item(10,6).Value = "Winter"
$ws.Cells.Item(10,7).Value = "Pah_WLow_depth.tif"
$ws.Cells.Item(10,8).Value = "Pah_WLow_velocity.tif"

# Row 11: scenario 9
$ws.Cells.Item(11,2).Value = 9
$ws.Cells.Item(11,3).Value = "Pahsimeroi"
$ws.Cells.Item(11,4).Value = "Chinook"
$ws.Cells.Item(11,5).Value = "Juvenile"
$ws.Cells.Item(11,6).Value = "Spring"
$ws.Cells.Item(11,7).Value = "Pah_1pt5_depth.tif"
$ws.Cells.Item(11,8).Value = "Pah_1pt5_velocity.tif"

# Row 12: scenario 10
$ws.Cells.Item(12,2).Value = 10
$ws.Cells.Item(12,3).Value = "Pahsimeroi"
$ws.Cells.Item(12,4).Value = "Chinook"
$ws.Cells.Item(12,5).Value = "Spawning"
$ws.Cells.Item(12,6).Value = "Summer"
$ws.Cells.Item(12,7).Value = "Pah_WLow_depth.tif"
$ws.Cells.Item(12,8).Value = "Pah_WLow_velocity.tif"

# Row 13: scenario 11
$ws.Cells.Item(13,2).Value = 11
$ws.Cells.Item(13,3).Value = "Pahsimeroi"
$ws.Cells.Item(13,4).Value = "Steelhead"
$ws.Cells.Item(13,5).Value = "Juvenile"
$ws.Cells.Item(13,6).Value = "Summer"
$ws.Cells.Item(13,7).Clear()
$ws.Cells.Item(13,8).Clear()

# Row 14: scenario 12
$ws.Cells.Item(14,2).Value = 12
$ws.Cells.Item(14,3).Value = "Pahsimeroi"
$ws.Cells.Item(14,4).Value = "Steelhead"
$ws.Cells.Item(14,5).Value = "Juvenile"
$ws.Cells.Item(14,6).Value = "Winter"
$ws.Cells.Item(14,7).Value = "Pah_WLow_depth.tif"
$ws.Cells.Item(14,8).Value = "Pah_WLow_velocity.tif"

# Row 15: scenario 13
$ws.Cells.Item(15,2).Value = 13
$ws.Cells.Item(15,3).Value = "Pahsimeroi"
$ws.Cells.Item(15,4).Value = "Steelhead"
$ws.Cells.Item(15,5).Value = "Juvenile"
$ws.Cells.Item(15,6).Value = "Spring"
$ws.Cells.Item(15,7).Value = "Pah_1pt5_depth.tif"
$ws.Cells.Item(15,8).Value = "Pah_1pt5_velocity.tif"

# Row 16: scenario 14
$ws.Cells.Item(16,2).Value = 14
$ws.Cells.Item(16,3).Value = "Pahsimeroi"
$ws.Cells.Item(16,4).Value = "Steelhead"
$ws.Cells.Item(16,5).Value = "Spawning"
$ws.Cells.Item(16,6).Value = "Spring"
$ws.Cells.Item(16,7).Value = "Pah_1pt5_depth.tif"
$ws.Cells.Item(16,8).Value = "Pah_1pt5_velocity.tif"

# Row 17: scenario 15
$ws.Cells.Item(17,2).Value = 15
$ws.Cells.Item(17,3).Value = "Upper Salmon"
$ws.Cells.Item(17,4).Value = "Chinook"
$ws.Cells.Item(17,5).Value = "Juvenile"
$ws.Cells.Item(17,6).Value = "Summer"
$ws.Cells.Item(17,7).Value = "US_Summer75_depth.tif"
$ws.Cells.Item(17,8).Value = "US_Summer75_velocity.tif"

# Row 18: scenario 16
$ws.Cells.Item(18,2).Value = 16
$ws.Cells.Item(18,3).Value = "Upper Salmon"
$ws.Cells.Item(18,4).Value = "Chinook"
$ws.Cells.Item(18,5).Value = "Juvenile"
$ws.Cells.Item(18,6).Value = "Winter"
$ws.Cells.Item(18,7).Value = "US_Winter75_depth.tif"
$ws.Cells.Item(18,8).Value = "US_Winter75_velocity.tif"

# Row 19: scenario 17
$ws.Cells.Item(19,2).Value = 17
$ws.Cells.Item(19,3).Value = "Upper Salmon"
$ws.Cells.Item(19,4).Value = "Chinook"
$ws.Cells.Item(19,5).Value = "Juvenile"
$ws.Cells.Item(19,6).Value = "Spring"
$ws.Cells.Item(19,7).Value = "US_1pt5year_depth.tif"
$ws.Cells.Item(19,8).Value = "US_1pt5year_velocity.tif"

# Row 20: scenario 18
$ws.Cells.Item(20,2).Value = 18
$ws.Cells.Item(20,3).Value = "Upper Salmon"
$ws.Cells.Item(20,4).Value = "Chinook"
$ws.Cells.Item(20,5).Value = "Spawning"
$ws.Cells.Item(20,6).Value = "Summer"
$ws.Cells.Item(20,7).Value = "US_Summer75_depth.tif"
$ws.Cells.Item(20,8).Value = "US_Summer75_velocity.tif"

# Row 21: scenario 19
$ws.Cells.Item(21,2).Value = 19
$ws.Cells.Item(21,3).Value = "Upper Salmon"
$ws.Cells.Item(21,4).Value = "Steelhead"
$ws.Cells.Item(21,5).Value = "Juvenile"
$ws.Cells.Item(21,6).Value = "Summer"
$ws.Cells.Item(21,7).Value = "US_Summer75_depth.tif"
$ws.Cells.Item(21,8).Value = "US_Summer75_velocity.tif"

# Row 22: scenario 20
$ws.Cells.Item(22,2).Value = 20
$ws.Cells.Item(22,3).Value = "Upper Salmon"
$ws.Cells.Item(22,4).Value = "Steelhead"
$ws.Cells.Item(22,5).Value = "Juvenile"
$ws.Cells.Item(22,6).Value = "Winter"
$ws.Cells.Item(22,7).Value = "US_Winter75_depth.tif"
$ws.Cells.Item(22,8).Value = "US_Winter75_velocity.tif"

# Row 23: scenario 21
$ws.Cells.Item(23,2).Value = 21
$ws.Cells.Item(23,3).Value = "Upper Salmon"
$ws.Cells.Item(23,4).Value = "Steelhead"
$ws.Cells.Item(23,5).Value = "Juvenile"
$ws.Cells.Item(23,6).Value = "Spring"
$ws.Cells.Item(23,7).Value = "US_1pt5year_depth.tif"
$ws.Cells.Item(23,8).Value = "US_1pt5year_velocity.tif"

# Row 24: scenario 22
$ws.Cells.Item(24,2).Value = 22
$ws.Cells.Item(24,3).Value = "Upper Salmon"
$ws.Cells.Item(24,4).Value = "Steelhead"
$ws.Cells.Item(24,5).Value = "Spawning"
$ws.Cells.Item(24,6).Value = "Spring"
$ws.Cells.Item(24,7).Value = "US_1pt5year_depth.tif"
$ws.Cells.Item(24,8).Value = "US_1pt5year_velocity.tif"

# ---------------------------------------------------------------
# Column I ("Date Completed") for rows 9-24: copy the date style
# from I15 onto every row, then set (or leave blank) the value.
# ---------------------------------------------------------------
$ws.Range("I15").Copy()
$ws.Range("I10:I24").PasteSpecial(-4122)
$ws.Cells.Item(10,9).Value = 43860
$ws.Cells.Item(11,9).Value = 43860
$ws.Cells.Item(12,9).Value = 43860
$ws.Cells.Item(14,9).Value = 43860
$ws.Cells.Item(15,9).Value = 43861
$ws.Cells.Item(16,9).Value = 43861
$ws.Cells.Item(17,9).Value = 43859
$ws.Cells.Item(18,9).Value = 43859
$ws.Cells.Item(19,9).Value = 43859
$ws.Cells.Item(20,9).Value = 43859
$ws.Cells.Item(21,9).Value = 43860
$ws.Cells.Item(22,9).Value = 43860
$ws.Cells.Item(23,9).Value = 43860
$ws.Cells.Item(24,9).Value = 43860
# I13 keeps the date style but no value (matches the source row).

# ---------------------------------------------------------------
# New blank, styled placeholder cells picked up by the edit:
# J10, J11, J14, J15 (centre style) and the new K10 cell, which
# introduces a left-aligned style.
# ---------------------------------------------------------------
$ws.Cells.Item(10,10).HorizontalAlignment = -4108
$ws.Cells.Item(11,10).HorizontalAlignment = -4108
$ws.Cells.Item(14,10).HorizontalAlignment = -4108
$ws.Cells.Item(15,10).HorizontalAlignment = -4108
$ws.Cells.Item(10,11).HorizontalAlignment = -4131

# ---------------------------------------------------------------
# Final selection, matching the saved cursor position.
# ---------------------------------------------------------------
$ws.Range("I9").Select()